$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B47 was mistakenly stored as text "4" -- fix it to be a real number
$ws.Range("B47").Value = 4

# Append a new annotation row (row 48) for Sunsi Wu
$ws.Range("A48").Value = "Sunsi Wu"

# politeness_score for this row is kept as literal text "1" (matches source data quirk)
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "1"
$ws.Range("B48").Style = "Normal"

$ws.Range("C48").Value = "so rude and misleading;willful misinterpretations"
$ws.Range("D48").Value = "CRT"
$ws.Range("E48").Value = "OTH"
$ws.Range("F48").Value = "f6e31c12-680e-4edf-b736-d4a426f6f32f"
$ws.Range("G48").Value = "Syg-YfWCW_annotated.xlsx"
$ws.Range("H48").Value = "It is incredible that the commenter continues to be so rude and misleading (should OpenReview have a moderating system?), and continues to frame this interaction as an attempt to convince *them* rather than to correct the constant series of willful misinterpretations and falsehoods that they manage to state about our work in every single interaction, in the hope that they do not mislead others."
